# Natmi following Dr Hou advice
# Update row 2 (target cluster corrected from FAPs to ECs, values recomputed),
# re-affirm row 3 with fresh FAPs-target values, and add a new row 4 for sCs target.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Vip"
$ws.Range("C2").Value = "Vipr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.119963
$ws.Range("H2").Value = 3.359889
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1382803333333333
$ws.Range("N2").Value = 0.414841
$ws.Range("O2").Value = 0.05497238930810543
$ws.Range("P2").Value = 0.05497238930810543
$ws.Range("Q2").Value = 0.154868856961
$ws.Range("R2").Value = 1.393819712649
$ws.Range("S2").Value = 0.05497238930810543
$ws.Range("T2").Value = 0.05497238930810543
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Vip"
$ws.Range("C3").Value = "Vipr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.119963
$ws.Range("H3").Value = 3.359889
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.696771
$ws.Range("N3").Value = 2.090313
$ws.Range("O3").Value = 0.2769964878394223
$ws.Range("P3").Value = 0.2769964878394223
$ws.Range("Q3").Value = 0.780357739473
$ws.Range("R3").Value = 7.023219655257
$ws.Range("S3").Value = 0.2769964878394223
$ws.Range("T3").Value = 0.2769964878394223
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Vip"
$ws.Range("C4").Value = "Vipr2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.119963
$ws.Range("H4").Value = 3.359889
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.680399333333333
$ws.Range("N4").Value = 5.041198
$ws.Range("O4").Value = 0.6680311228524723
$ws.Range("P4").Value = 0.6680311228524722
$ws.Range("Q4").Value = 1.881985078558
$ws.Range("R4").Value = 16.937865707022
$ws.Range("S4").Value = 0.6680311228524723
$ws.Range("T4").Value = 0.6680311228524722

$wb.Save()
